# Auto-generated edit script applying scheduled runner updates to Leve profit sheets
$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H2").Value = 179.25
$ws.Range("I2").Value = 179.25
$ws.Range("J2").Value = 0
$ws.Range("K2").Value = 179.25
$ws.Range("L2").Value = 0
$ws.Range("M2").Value = -66.25
$ws.Range("N2").Value = $null

$ws.Range("H40").Value = 3347.2222
$ws.Range("I40").Value = 1187.5
$ws.Range("J40").Value = 5075
$ws.Range("K40").Value = 1187.5
$ws.Range("L40").Value = 5075
$ws.Range("M40").Value = -1012.5
$ws.Range("N40").Value = -5425

$ws.Range("H43").Value = 20000
$ws.Range("I43").Value = 20000
$ws.Range("J43").Value = 0
$ws.Range("K43").Value = 20000
$ws.Range("L43").Value = 0
$ws.Range("M43").Value = -19931
$ws.Range("N43").Value = $null

$ws.Range("H53").Value = 6556.65
$ws.Range("I53").Value = 1170.9166
$ws.Range("K53").Value = 1170.9166
$ws.Range("M53").Value = -533.9166

$ws.Range("H113").Value = 9693.75
$ws.Range("I113").Value = 10573.091
$ws.Range("J113").Value = 8619
$ws.Range("K113").Value = 10573.091
$ws.Range("L113").Value = 8619
$ws.Range("M113").Value = -7319.091
$ws.Range("N113").Value = -15127

$ws.Range("H115").Value = 474.44446
$ws.Range("I115").Value = 452.625
$ws.Range("J115").Value = 649
$ws.Range("K115").Value = 1357.875
$ws.Range("L115").Value = 1947
$ws.Range("M115").Value = 209.125
$ws.Range("N115").Value = -5081

$ws.Range("H116").Value = 4241.263
$ws.Range("I116").Value = 4458.4165
$ws.Range("J116").Value = 3869
$ws.Range("K116").Value = 4458.4165
$ws.Range("L116").Value = 3869
$ws.Range("M116").Value = -1016.4165
$ws.Range("N116").Value = -10753

$ws.Range("H137").Value = 20061.777
$ws.Range("J137").Value = 21913.125
$ws.Range("L137").Value = 65739.375
$ws.Range("N137").Value = -70839.375

$ws.Range("H138").Value = 4323.5884
$ws.Range("J138").Value = 2513.2727
$ws.Range("L138").Value = 7539.8181
$ws.Range("N138").Value = -17819.8181

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H63").Value = 2698.25
$ws.Range("I63").Value = 2264.3333
$ws.Range("J63").Value = 4000
$ws.Range("K63").Value = 2264.3333
$ws.Range("L63").Value = 4000
$ws.Range("M63").Value = -1578.3333
$ws.Range("N63").Value = -5372

$ws.Range("H66").Value = 2698.25
$ws.Range("I66").Value = 2264.3333
$ws.Range("J66").Value = 4000
$ws.Range("K66").Value = 11321.6665
$ws.Range("L66").Value = 20000
$ws.Range("M66").Value = -7889.666499999999
$ws.Range("N66").Value = -26864

$ws.Range("H74").Value = 14971.7
$ws.Range("I74").Value = 899.26666
$ws.Range("K74").Value = 899.26666
$ws.Range("M74").Value = -25.26666

$ws.Range("H77").Value = 14971.7
$ws.Range("I77").Value = 899.26666
$ws.Range("K77").Value = 4496.3333
$ws.Range("M77").Value = -128.3333000000002

$ws.Range("H110").Value = 2722.074
$ws.Range("I110").Value = 726.0952
$ws.Range("K110").Value = 726.0952
$ws.Range("M110").Value = 1318.9048

$ws.Range("H132").Value = 3348186.2
$ws.Range("I132").Value = 5148.0586
$ws.Range("K132").Value = 15444.1758
$ws.Range("M132").Value = -12914.1758

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H14").Value = 1428.7273
$ws.Range("I14").Value = 679.5
$ws.Range("J14").Value = 1856.8572
$ws.Range("K14").Value = 679.5
$ws.Range("L14").Value = 1856.8572
$ws.Range("M14").Value = -507.5
$ws.Range("N14").Value = -2200.8572

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H16").Value = 4339.353
$ws.Range("J16").Value = 9962.333000000001
$ws.Range("L16").Value = 9962.333000000001
$ws.Range("N16").Value = -10536.333

$ws.Range("H17").Value = 29999.75
$ws.Range("I17").Value = 50999.5
$ws.Range("K17").Value = 50999.5
$ws.Range("M17").Value = -50825.5

$ws.Range("H25").Value = 1800
$ws.Range("I25").Value = 1800
$ws.Range("K25").Value = 1800
$ws.Range("M25").Value = -1626

$ws.Range("H113").Value = 4339.353
$ws.Range("J113").Value = 9962.333000000001
$ws.Range("L113").Value = 9962.333000000001
$ws.Range("N113").Value = -14302.333

$ws.Range("H132").Value = 8792.1
$ws.Range("I132").Value = 2416.2856
$ws.Range("K132").Value = 7248.8568
$ws.Range("M132").Value = -4718.8568

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H4").Value = 73366920
$ws.Range("I4").Value = 93285150
$ws.Range("K4").Value = 279855450
$ws.Range("M4").Value = -279855338

$ws.Range("H122").Value = 13454206
$ws.Range("J122").Value = 4055187.2
$ws.Range("L122").Value = 36496684.8
$ws.Range("N122").Value = -36501584.8

$ws.Range("H131").Value = 1480.79
$ws.Range("I131").Value = 1050
$ws.Range("J131").Value = 1494.1134
$ws.Range("K131").Value = 3150
$ws.Range("L131").Value = 4482.3402
$ws.Range("M131").Value = 1890
$ws.Range("N131").Value = -14562.3402

$ws.Range("H139").Value = 16539.637
$ws.Range("I139").Value = 19015.223
$ws.Range("K139").Value = 57045.66900000001
$ws.Range("M139").Value = -51905.66900000001

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H80").Value = 14673.2
$ws.Range("I80").Value = 7109
$ws.Range("J80").Value = 35474.75
$ws.Range("K80").Value = 7109
$ws.Range("L80").Value = 35474.75
$ws.Range("M80").Value = -6111
$ws.Range("N80").Value = -37470.75

$ws.Range("H83").Value = 14673.2
$ws.Range("I83").Value = 7109
$ws.Range("J83").Value = 35474.75
$ws.Range("K83").Value = 35545
$ws.Range("L83").Value = 177373.75
$ws.Range("M83").Value = -30553
$ws.Range("N83").Value = -187357.75

$ws.Range("H102").Value = 5517.885
$ws.Range("I102").Value = 4453.4736
$ws.Range("J102").Value = 8407
$ws.Range("K102").Value = 4453.4736
$ws.Range("L102").Value = 8407
$ws.Range("M102").Value = -2831.4736
$ws.Range("N102").Value = -11651

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H20").Value = 3761
$ws.Range("I20").Value = 1812.25
$ws.Range("J20").Value = 5320
$ws.Range("K20").Value = 1812.25
$ws.Range("L20").Value = 5320
$ws.Range("M20").Value = -1586.25
$ws.Range("N20").Value = -5772

$ws.Range("H26").Value = 0
$ws.Range("I26").Value = 0
$ws.Range("K26").Value = 0
$ws.Range("M26").Value = $null

$ws.Range("H122").Value = 7245.091
$ws.Range("I122").Value = 5499.3
$ws.Range("J122").Value = 8699.916999999999
$ws.Range("K122").Value = 16497.9
$ws.Range("L122").Value = 26099.751
$ws.Range("M122").Value = -14047.9
$ws.Range("N122").Value = -30999.751

$ws.Range("H132").Value = 1151218
$ws.Range("I132").Value = 2165.0356
$ws.Range("J132").Value = 5747430
$ws.Range("K132").Value = 6495.1068
$ws.Range("L132").Value = 17242290
$ws.Range("M132").Value = -3965.1068
$ws.Range("N132").Value = -17247350

$ws.Range("H136").Value = 13393.415
$ws.Range("J136").Value = 11140.786
$ws.Range("L136").Value = 33422.358
$ws.Range("N136").Value = -38522.358
